$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CSV "Cargo" reference sheet renames two cargo names:
#   战斗机 (fighter)  -> 攻击机 (attack craft)
#   轰炸机 (bomber)   -> 防御机 (defense craft)
$ws.Range("B4").Value = "攻击机"
$ws.Range("B5").Value = "防御机"

# Move the active selection to B6, matching where the user ended up editing
$ws.Range("B6").Select()
